# Adds the "Vacate a default judgment of foreclosure" easy-form entry to
# the lookup table on Sheet1, just above the existing
# "Vacate a default judgment within 30 days" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 50, pushing the current rows 50-51
# (and their formatting / hyperlink style) down to rows 51-52.
$ws.Rows(50).Insert()

# Populate the new row with the name/url pair for the new form.
$ws.Range("A50").Value = "Vacate a default judgment of foreclosure"
$ws.Range("B50").Value = "https://easyforms.illinoislegalaid.org/start/VacateDefaultForeclosure/vacate_default_foreclosure"

# The row insert does not renumber the worksheet's <hyperlinks> collection,
# so rebuild it from scratch against the new row layout: the two hyperlinks
# that used to target B50/B51 now belong on B51/B52, every other hyperlink
# keeps its original target cell, and the new form gets a hyperlink on B50.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.illinoislegalaid.org/legal-information/appearance") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B18"), "https://www.illinoislegalaid.org/legal-information/fee-waiver") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B15"), "https://www.illinoislegalaid.org/legal-information/end-illegal-lockout-demand") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B43"), "https://www.illinoislegalaid.org/legal-information/security-deposit-demand-letter") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B20"), "https://www.illinoislegalaid.org/legal-information/housing-discrimination-complaint-idhr") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B47"), "https://www.illinoislegalaid.org/legal-information/stop-wage-assignment-letter") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B36"), "https://www.illinoislegalaid.org/legal-information/request-time-work-due-domestic-abuse-letter") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-circuit-court") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-appellate-court") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B13"), "https://www.illinoislegalaid.org/legal-information/e-filing-exemption-supreme-court") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B39"), "https://www.illinoislegalaid.org/legal-information/respond-lawsuit") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B52"), "https://www.illinoislegalaid.org/legal-information/voluntary-acknowledgment-parentage-vap") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B21"), "https://www.illinoislegalaid.org/legal-information/interpreter-request") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B23"), "https://www.illinoislegalaid.org/legal-information/motion") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B48"), "https://www.illinoislegalaid.org/legal-information/transfer-death-instrument-or-todi") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B32"), "https://www.illinoislegalaid.org/legal-information/power-attorney-agent-resign-letter") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B33"), "https://www.illinoislegalaid.org/legal-information/power-attorney-revocation") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B31"), "https://www.illinoislegalaid.org/legal-information/power-attorney-property") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B30"), "https://www.illinoislegalaid.org/legal-information/power-attorney-health-care") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B27"), "https://www.illinoislegalaid.org/legal-information/order-protection") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B25"), "https://www.illinoislegalaid.org/legal-information/name-change-adult") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.illinoislegalaid.org/legal-information/cannabis-expungement") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B14"), "https://www.illinoislegalaid.org/legal-information/emergency-order-protection-cook-county") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B44"), "https://www.illinoislegalaid.org/legal-information/short-term-guardian-appointment") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B35"), "https://www.illinoislegalaid.org/legal-information/remove-eviction-public-record") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B41"), "https://www.illinoislegalaid.org/legal-information/respond-eviction") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B45"), "https://www.illinoislegalaid.org/legal-information/small-claims-complaint") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B42"), "https://www.illinoislegalaid.org/legal-information/security-deposit-complaint") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://www.illinoislegalaid.org/legal-information/criminal-court-fee-waiver") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B19"), "https://www.illinoislegalaid.org/legal-information/financial-affidavit") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B37"), "https://www.illinoislegalaid.org/legal-information/special-process-server-request") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B51"), "https://www.illinoislegalaid.org/legal-information/vacate-default-judgment-within-30-days") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B24"), "https://www.illinoislegalaid.org/legal-information/motion-continue-or-extend-time") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://www.illinoislegalaid.org/legal-information/debt-collector-letter") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B49"), "https://www.illinoislegalaid.org/legal-information/unemployment-benefits-request-hearing") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B38"), "https://www.illinoislegalaid.org/legal-information/child-support-modification") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B50"), "https://easyforms.illinoislegalaid.org/start/VacateDefaultForeclosure/vacate_default_foreclosure") | Out-Null

# Match the author's final cursor position/selection.
$ws.Range("B56").Select()
